# Morgonmote.docx -- sync fran lordagen + morgonmote.
#
# Paragraph 2 ("Idag sa ska vi titta...") is retyped into three runs (the
# "lo|gin" split also fixes the old "loggin" typo -> "login").
# Paragraph 3's date "2013-11-25" is corrected to "2016-11-25" (split into
# "2016" / "-11-25" runs).
# Paragraph 4 ("Igar sa borjade...") is retyped with proofErr spell-check
# wrappers around "griden"/"gridden" (text unchanged).
# Six brand-new paragraphs (two more day-entries) are appended after it,
# ending with the paragraph that keeps the original _GoBack bookmark.

$d = $word.ActiveDocument

if ($d.Paragraphs.Count -lt 4) {
    throw "unexpected document shape: expected at least 4 paragraphs, found $($d.Paragraphs.Count)"
}

$p2Text = $d.Paragraphs(2).Range.Text
$p4Text = $d.Paragraphs(4).Range.Text
if ($p2Text -notmatch "Idag s.* ska vi titta") {
    throw "paragraph 2 does not look like the expected 'Idag sa ska vi titta...' paragraph"
}
if ($p4Text -notmatch "griden") {
    throw "paragraph 4 does not look like the expected 'Igar sa borjade...' paragraph"
}

# Replace paragraphs 2 through 4 (inclusive) in one shot with their new
# contents plus the newly authored day entries that follow them.
$startPos = $d.Paragraphs(2).Range.Start
$endPos = $d.Paragraphs(4).Range.End
$target = $d.Range($startPos, $endPos)

$newXml = @'
<w:p><w:r><w:t xml:space="preserve">Idag så ska vi titta på film med Terese på morgonen, som kommer handla om ingenjörer. På eftermiddagen så kommer Emma jobba med att </w:t></w:r><w:r><w:t>få i gång registreringen och lo</w:t></w:r><w:r><w:t>gin. Dennis ska fixa den generella layouten för fönstret. Och Pontus latar sig på Dreamhack.</w:t></w:r></w:p><w:p><w:r><w:t>2016</w:t></w:r><w:r><w:t>-11-25</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Igår så började vi programmera i java Dennis började med att skapa fönstret och den generella </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>griden</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, och Emma började med registreringen. Idag ska Emma se till så att formuläret fungerar. Dennis ska fina till den generella </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gridden</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> och möjligt vi börja med att skapa händelser.</w:t></w:r></w:p><w:p><w:r><w:t>2016-11-26</w:t></w:r></w:p><w:p><w:r><w:t>Igår så gjordes registrering grundklart, dvs det går att registrera sig men finns ingen felhantering av det än. Login påbörjades. Vi fixade databas-anslutningen. Dennis påbörjade att lägga till händelser.</w:t></w:r></w:p><w:p><w:r><w:t>Idag, lördag ska vi visa upp oss på öppet hus och helt enkelt bara fortsätta arbeta. Händelser och login.</w:t></w:r></w:p><w:p><w:r><w:t>2016-11-28</w:t></w:r></w:p><w:p><w:r><w:t>Öppet hus på lördagen var välbesökt, dessvärre var det så välbesökt att vi inte hann göra något arbete.</w:t></w:r></w:p><w:p><w:r><w:t>Vi får fortsätta med det vi skulle gjort i lördags. Pontus ska göra knappar till vyn, men det ska inte fungera än.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$target.InsertXML($newXml)
